# Ben Dudley timesheet update (25.01.16)
# - Time Taken for "Implement Junit tests from specification" row: 1:30 -> 4:30
# - Time Taken for "Testing for all targeted platforms and bug fixes" row: 6:00 -> 18:00
# - Updated the "Total hours accounted this week" summary text
# - Selection cursor left on G2 (instead of I3) when the sheet was last saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0.1875
$ws.Range("F4").Value = 0.75

$ws.Range("G1").Value = "Total hours accounted this week: 29.30"

[void]$ws.Range("G2").Select()
